$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.040.78"
$ws.Range("E2").Value = "  -2.30%  "

$ws.Range("D3").Value = "3.565.61"

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "620.06"
$ws.Range("E5").Value = "  -7.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.27"
$ws.Range("E6").Value = "  -3.61%  "

$ws.Range("D7").Value = "3.562.45"
$ws.Range("E7").Value = "  -3.36%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  -2.29%  "

$ws.Range("E10").Value = "  -3.08%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.97"
$ws.Range("E11").Value = "  -2.00%  "

$ws.Range("E12").Value = "  -1.73%  "

$ws.Range("E13").Value = "  -3.45%  "

$ws.Range("D14").Value = "4.168.49"
$ws.Range("E14").Value = "  -3.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "32.13"
$ws.Range("E15").Value = "  -2.05%  "

$ws.Range("D16").Value = "3.557.84"
$ws.Range("E16").Value = "  -4.39%  "

$ws.Range("D17").Value = "68.069.05"
$ws.Range("E17").Value = "  -2.28%  "

$ws.Range("E18").Value = "  -1.01%  "

$ws.Range("E19").Value = "  -0.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.64"
$ws.Range("E20").Value = "  -3.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "459.30"
$ws.Range("E21").Value = "  -2.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.72"
$ws.Range("E22").Value = "  -0.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.646"
$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.80"
$ws.Range("E24").Value = "  -2.62%  "

$ws.Range("D25").Value = "3.707.77"
$ws.Range("E25").Value = "  -3.36%  "

$ws.Range("E26").Value = "  +0.14%  "

$ws.Range("E27").Value = "  -2.02%  "

$ws.Range("E28").Value = "  -8.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.38"
$ws.Range("E29").Value = "  -7.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.57"
$ws.Range("E30").Value = "  -3.84%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.64"
$ws.Range("E31").Value = "  -3.57%  "

$ws.Range("E32").Value = "  +0.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.02"
$ws.Range("E33").Value = "  -2.77%  "

$ws.Range("E34").Value = "  -4.36%  "

$ws.Range("E35").Value = "  -4.38%  "

$ws.Range("D36").Value = "3.565.05"
$ws.Range("E36").Value = "  -3.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.19"
$ws.Range("E37").Value = "  -4.39%  "

$ws.Range("E38").Value = "  -4.08%  "

$ws.Range("E39").Value = "  -0.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "178.22"
$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("E41").Value = "  -0.10%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0889"
$ws.Range("E42").Value = "  -1.67%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.62"
$ws.Range("E43").Value = "  -7.82%  "

$ws.Range("E44").Value = "  -5.96%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.897"
$ws.Range("E45").Value = "  -3.98%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.83"
$ws.Range("E46").Value = "  +4.61%  "

$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.04"
$ws.Range("E47").Value = "  -2.09%  "

$ws.Range("E48").Value = "  -6.73%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.73"
$ws.Range("E49").Value = "  -1.48%  "

$ws.Range("E50").Value = "  -5.74%  "

$ws.Range("E51").Value = "  -5.35%  "
